$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect numeric-looking Price cells as Text first, so assigning the new
# literal digit string does not get auto-converted to a Number by Excel
# (which would lose trailing zeros / introduce float rounding / sci-notation).
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'

$ws.Range('D2').Value = '25.787.17'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').Value = '1.635.31'
$ws.Range('E3').Value = '  -0.16%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '215.25'
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('D6').Value = '0.505'
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').Value = '0.0641'
$ws.Range('E9').Value = '  -0.20%  '
$ws.Range('D10').Value = '19.83'
$ws.Range('E10').Value = '  +0.51%  '
$ws.Range('E11').Value = '  +0.05%  '
$ws.Range('E12').Value = '  -0.79%  '
$ws.Range('D13').Value = '1.636.22'
$ws.Range('E13').Value = '  -0.06%  '
$ws.Range('D14').Value = '1.859.88'
$ws.Range('E14').Value = '  -0.26%  '
$ws.Range('D15').Value = '0.556'
$ws.Range('E15').Value = '  -0.66%  '
$ws.Range('D16').Value = '0.0₃0776'
$ws.Range('E16').Value = '  +1.80%  '
$ws.Range('D17').Value = '63.16'
$ws.Range('E17').Value = '  +0.29%  '
$ws.Range('D18').Value = '25.793.69'
$ws.Range('E18').Value = '  -0.33%  '
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('D20').Value = '4.44'
$ws.Range('E20').Value = '  +2.65%  '
$ws.Range('D21').Value = '194.40'
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('D22').Value = '9.96'
$ws.Range('E22').Value = '  +0.88%  '
$ws.Range('D23').Value = '6.15'
$ws.Range('E23').Value = '  +0.89%  '
$ws.Range('D24').Value = '1.01'
$ws.Range('E24').Value = '  +0.25%  '
$ws.Range('E25').Value = '  -0.80%  '
$ws.Range('D26').Value = '139.90'
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('D27').Value = '0.121'
$ws.Range('E27').Value = '  -3.94%  '
$ws.Range('D28').Value = '6.85'
$ws.Range('E28').Value = '  +0.44%  '
$ws.Range('D29').Value = '15.60'
$ws.Range('E29').Value = '  +1.31%  '
$ws.Range('D31').Value = '0.0492'
$ws.Range('E31').Value = '  +0.70%  '
$ws.Range('D32').Value = '3.35'
$ws.Range('E32').Value = '  +1.59%  '
$ws.Range('D33').Value = '3.27'
$ws.Range('E33').Value = '  +1.47%  '
$ws.Range('E34').Value = '  +1.49%  '
$ws.Range('D36').Value = '0.898'
$ws.Range('E36').Value = '  -0.58%  '
$ws.Range('D37').Value = '2.58'
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('E38').Value = '  +0.25%  '
$ws.Range('D39').Value = '1.108.17'
$ws.Range('E39').Value = '  -1.62%  '
$ws.Range('D40').Value = '0.0157'
$ws.Range('E40').Value = '  +0.33%  '
$ws.Range('E41').Value = '  +0.43%  '
$ws.Range('E42').Value = '  +0.33%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '0.802'
$ws.Range('E43').Value = '  +0.26%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value = '99.24'
$ws.Range('E44').Value = '  +1.17%  '
$ws.Range('D45').Value = '0.0₆0108'
$ws.Range('E45').Value = '  -4.37%  '
$ws.Range('D46').Value = '55.19'
$ws.Range('E46').Value = '  -0.36%  '
$ws.Range('D47').Value = '2.49'
$ws.Range('E47').Value = '  +12.70%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = '0.418'
$ws.Range('E48').Value = '  -2.08%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '0.0504'
$ws.Range('E49').Value = '  -0.08%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '7.66'
$ws.Range('E50').Value = '  -0.25%  '
$ws.Range('E51').Value = '  -0.76%  '
